# Add federal-holiday ("X") marking to the FRI column (R:S) on the
# "16-End" sheet, mirroring the styling already used for the adjacent
# SAT column (T:U), plus narrowing the FRI columns to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("16-End")

# Narrow columns R (18) and S (19) from 4.5 to 2.5 characters, matching
# the width already used by the neighboring weekend columns (e.g. col T/20).
$refWidth = $ws.Columns.Item(20).ColumnWidth
$ws.Columns.Item(18).ColumnWidth = $refWidth
$ws.Columns.Item(19).ColumnWidth = $refWidth

# Rows whose FRI (R:S) cells should be filled with "X" (federal holiday),
# matching the pattern already used on the SAT (T:U) columns in those rows.
$xRows = @(5,6,8,9,11,12,14,15,17,18,20,21,23,24,26,27)

for ($row = 2; $row -le 27; $row++) {
    $srcRange = $ws.Range("T" + $row + ":U" + $row)
    $dstRange = $ws.Range("R" + $row + ":S" + $row)

    # Copy the SAT column's formatting (fill, borders, font) onto the
    # FRI column for this row, leaving existing values untouched.
    $srcRange.Copy()
    $dstRange.PasteSpecial(-4122)

    if ($xRows -contains $row) {
        $dstRange.Value = "X"
    }
}

$excel.CutCopyMode = 0
